$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A128").Value = 0.007403190136094686
$ws.Range("E128").Value = 0.007403190136094686
$ws.Range("A129").Value = 0.007478349426816459
$ws.Range("E129").Value = 0.01488153956291114
$ws.Range("A130").Value = 0.01121752414022469
$ws.Range("E130").Value = 0.02609906370313583
$ws.Range("A131").Value = 0.02243504828044938
$ws.Range("E131").Value = 0.04853411198358521
$ws.Range("A132").Value = 0.04487009656089876
$ws.Range("E132").Value = 0.09340420854448397
$ws.Range("A133").Value = 0.0859972604438532
$ws.Range("E133").Value = 0.1794014689883372
$ws.Range("A134").Value = 0.145827813822921
$ws.Range("E134").Value = 0.3252292828112581
$ws.Range("A135").Value = 0.2213196844061383
$ws.Range("E135").Value = 0.5465489672173963
$ws.Range("A136").Value = 0.2223362138131503
$ws.Range("E136").Value = 0.7688851810305466
$ws.Range("A137").Value = 0.1061136446055359
$ws.Range("E137").Value = 0.8749988256360824
$ws.Range("A138").Value = 0.04981557789039145
$ws.Range("E138").Value = 0.9248144035264739
$ws.Range("A139").Value = 0.03672282944665851
$ws.Range("E139").Value = 0.9615372329731324
$ws.Range("A140").Value = 0.02040950539549758
$ws.Range("E140").Value = 0.98194673836863
$ws.Range("A141").Value = 0.007006724877537331
$ws.Range("E141").Value = 0.9889534632461674
$ws.Range("A142").Value = 0.002529110132787677
$ws.Range("E142").Value = 0.9914825733789551
$ws.Range("A143").Value = 0.005997711399597522
$ws.Range("E143").Value = 0.9974802847785526
$ws.Range("A144").Value = 0.00004885353896915275
$ws.Range("E144").Value = 0.9975291383175217
$ws.Range("A145").Value = 0.0004791404783513058
$ws.Range("E145").Value = 0.998008278795873
$ws.Range("A146").Value = 0.0003363378259799362
$ws.Range("E146").Value = 0.9983446166218529
$ws.Range("A147").Value = 0.00165538337814706
$ws.Range("A153").Value = 0.02860161267339187
$ws.Range("E153").Value = 0.02860161267339187
$ws.Range("A154").Value = 0.04916934053352857
$ws.Range("E154").Value = 0.07777095320692044
$ws.Range("A155").Value = 0.05443976551866465
$ws.Range("E155").Value = 0.1322107187255851
$ws.Range("A156").Value = 0.07064217027904851
$ws.Range("E156").Value = 0.2028528890046336
$ws.Range("A157").Value = 0.1264105978487861
$ws.Range("E157").Value = 0.3292634868534197
$ws.Range("A158").Value = 0.2196036616825461
$ws.Range("E158").Value = 0.5488671485359657
$ws.Range("A159").Value = 0.138486484892857
$ws.Range("E159").Value = 0.6873536334288227
$ws.Range("A160").Value = 0.07486106538648636
$ws.Range("E160").Value = 0.7622146988153091
$ws.Range("A161").Value = 0.06061824060379435
$ws.Range("E161").Value = 0.8228329394191034
$ws.Range("A162").Value = 0.06406273800690865
$ws.Range("E162").Value = 0.8868956774260121
$ws.Range("A163").Value = 0.03310353638342514
$ws.Range("E163").Value = 0.9199992138094373
$ws.Range("A164").Value = 0.0329806941079931
$ws.Range("E164").Value = 0.9529799079174304
$ws.Range("A165").Value = 0.01699743996698
$ws.Range("E165").Value = 0.9699773478844104
$ws.Range("A166").Value = 0.01710652390756364
$ws.Range("E166").Value = 0.987083871791974
$ws.Range("A167").Value = 0.008827937281647855
$ws.Range("E167").Value = 0.9959118090736219
$ws.Range("A168").Value = 0.001907494852908659
$ws.Range("E168").Value = 0.9978193039265305
$ws.Range("A169").Value = 0.001365023364600787
$ws.Range("E169").Value = 0.9991843272911313
$ws.Range("A170").Value = 0.0002476500272709851
$ws.Range("E170").Value = 0.9994319773184023
$ws.Range("E171").Value = 0.9994319773184023
$ws.Range("A172").Value = 0.00000196547640691258
$ws.Range("E172").Value = 0.9994339427948091
$ws.Range("A173").Value = 0.0005660572051908232
$ws.Range("A179").Value = 0.02778737170141132
$ws.Range("E179").Value = 0.02778737170141132
$ws.Range("A180").Value = 0.04761911039281185
$ws.Range("E180").Value = 0.07540648209422317
$ws.Range("A181").Value = 0.05272337071008804
$ws.Range("E181").Value = 0.1281298528043112
$ws.Range("A182").Value = 0.06841580738799497
$ws.Range("E182").Value = 0.1965456601923062
$ws.Range("A183").Value = 0.1224264482943176
$ws.Range("E183").Value = 0.3189721084866237
$ws.Range("A184").Value = 0.2126816876751624
$ws.Range("E184").Value = 0.5316537961617861
$ws.Range("A185").Value = 0.1341212947314859
$ws.Range("E185").Value = 0.6657750908932719
$ws.Range("A186").Value = 0.07250106119048436
$ws.Range("E186").Value = 0.7382761520837563
$ws.Range("A187").Value = 0.05870756226749382
$ws.Range("E187").Value = 0.7969837143512502
$ws.Range("A188").Value = 0.06204339148569305
$ws.Range("E188").Value = 0.8590271058369432
$ws.Range("A189").Value = 0.03197281243161583
$ws.Range("E189").Value = 0.890999918268559
$ws.Range("A190").Value = 0.04493581445380987
$ws.Range("E190").Value = 0.9359357327223689
$ws.Range("A191").Value = 0.02315899929078201
$ws.Range("E191").Value = 0.9590947320131509
$ws.Range("A192").Value = 0.02330730230878138
$ws.Range("E192").Value = 0.9824020343219322
$ws.Range("A193").Value = 0.01202770432201133
$ws.Range("E193").Value = 0.9944297386439436
$ws.Range("A194").Value = 0.00259958712439789
$ws.Range("E194").Value = 0.9970293257683415
$ws.Range("A195").Value = 0.001860049407974352
$ws.Range("E195").Value = 0.9988893751763158
$ws.Range("A196").Value = 0.0003374717565141277
$ws.Range("E196").Value = 0.9992268469328299
$ws.Range("E197").Value = 0.9992268469328299
$ws.Range("A198").Value = 0.000001977373573324967
$ws.Range("E198").Value = 0.9992288243064033
$ws.Range("A199").Value = 0.000771175693596737
